$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (not be coerced to a
# number). Force text entry with a leading apostrophe, then restore the
# original cell formatting (which the apostrophe/number-format dance
# otherwise perturbs) by re-pasting formats from a same-style neighbor.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("D2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 17.05.2024"

# Row 6
$ws.Range("B6").Value = "18.05."
$ws.Range("C6").Value = "19.05."
$ws.Range("D6").Value = "PAYPAL DLPTPO"
$ws.Range("E6").Value = "35,67-"

# Row 7
$ws.Range("B7").Value = "20.05."
$ws.Range("C7").Value = "21.05."
$ws.Range("D7").Value = "KARTENZ./20.05 REWE RO"
$ws.Range("E7").Value = "93,57-"

# Row 8
$ws.Range("B8").Value = "23.05."
$ws.Range("C8").Value = "24.05."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 28749976"
$ws.Range("E8").Value = "84,23-"

# Row 9 - newly populated transaction row, copy style from row 8 so
# formatting (font/number format/borders) matches the other data rows
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)

$ws.Range("B9").Value = "24.05."
$ws.Range("C9").Value = "25.05."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,72-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 27.05.2024"
$ws.Range("E12").Value = "238,19-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.06.2024"
